# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
# Both sheets carry duplicate data, and the same rows/values changed in each.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 1776
    4  = 1678
    6  = 235
    9  = 44
    11 = 20
    12 = 83
    14 = 239
    16 = 33
    17 = 23
    18 = 70
    19 = 204
    20 = 33
    21 = 435
    22 = 329
    25 = 27
    27 = 711
    28 = 2504
    31 = 503
    32 = 818
    34 = 438
    35 = 247
    37 = 425
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
